$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2022 events: Korea and China swap places (Korea now ranks above China)
$ws.Range("A2").Value = "Korea"
$ws.Range("B2").Value = 46
$ws.Range("A3").Value = "China"
$ws.Range("B3").Value = 41

# Updated scores for the remaining regions that kept their rank
$ws.Range("B4").Value = 29
$ws.Range("B5").Value = 29
$ws.Range("B6").Value = 20

# Drop the regions that did not take part in the 2022 events
$ws.Rows("13:16").Delete()
